$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Text fix: the standalone "Krakow (Cracow)" label becomes "Krakow" ---
# (the "Krakow (Cracow), Poland" entry is untouched)
$ws.Range("D226").Value = "Krakow"

# --- 2. Scroll / selection housekeeping so the view lands further down the
#        (now longer) sheet, matching the other same-length yearly datasets ---
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 209
$win.ScrollColumn = 1
$ws.Range("A226").Select() | Out-Null

# --- 3. Slightly narrower default column width across the sheet ---
$ws.Columns("A:AMK").ColumnWidth = 13.333333333333334
